# Update NATMI LR-pairs values with new TPM-derived values
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 3.062550666666667
$ws.Range("H2").Value = 9.187652
$ws.Range("I2").Value = 0.06849600470812313
$ws.Range("J2").Value = 0.06849600470812313
$ws.Range("M2").Value = 1.646588666666666
$ws.Range("N2").Value = 4.939766
$ws.Range("O2").Value = 0.039310317935267
$ws.Range("P2").Value = 0.039310317935267
$ws.Range("Q2").Value = 5.042761218825778
$ws.Range("R2").Value = 45.38485096943199
$ws.Range("S2").Value = 0.002692599722371865
$ws.Range("T2").Value = 0.002692599722371865
$ws.Range("G3").Value = 3.062550666666667
$ws.Range("H3").Value = 9.187652
$ws.Range("I3").Value = 0.06849600470812313
$ws.Range("J3").Value = 0.06849600470812313
$ws.Range("O3").Value = 0.278787195370394
$ws.Range("P3").Value = 0.278787195370394
$ws.Range("Q3").Value = 35.76305995372712
$ws.Range("R3").Value = 321.867539583544
$ws.Range("S3").Value = 0.01909580904665495
$ws.Range("T3").Value = 0.01909580904665495
$ws.Range("G4").Value = 3.062550666666667
$ws.Range("H4").Value = 9.187652
$ws.Range("I4").Value = 0.06849600470812313
$ws.Range("J4").Value = 0.06849600470812313
$ws.Range("M4").Value = 0.7553226666666667
$ws.Range("N4").Value = 2.265968
$ws.Range("O4").Value = 0.01803241742850595
$ws.Range("P4").Value = 0.01803241742850595
$ws.Range("Q4").Value = 2.313213936348445
$ws.Range("R4").Value = 20.818925427136
$ws.Range("S4").Value = 0.001235148549081785
$ws.Range("T4").Value = 0.001235148549081785
$ws.Range("G5").Value = 3.062550666666667
$ws.Range("H5").Value = 9.187652
$ws.Range("I5").Value = 0.06849600470812313
$ws.Range("J5").Value = 0.06849600470812313
$ws.Range("M5").Value = 27.21325766666666
$ws.Range("N5").Value = 81.63977299999999
$ws.Range("O5").Value = 0.6496836961088899
$ws.Range("P5").Value = 0.6496836961088899
$ws.Range("Q5").Value = 83.34198040922176
$ws.Range("R5").Value = 750.0778236829959
$ws.Range("S5").Value = 0.04450073750746536
$ws.Range("T5").Value = 0.04450073750746536
$ws.Range("G6").Value = 3.062550666666667
$ws.Range("H6").Value = 9.187652
$ws.Range("I6").Value = 0.06849600470812313
$ws.Range("J6").Value = 0.06849600470812313
$ws.Range("M6").Value = 0.5942236666666667
$ws.Range("N6").Value = 1.782671
$ws.Range("O6").Value = 0.01418637315694314
$ws.Range("P6").Value = 0.01418637315694314
$ws.Range("Q6").Value = 1.819840086499111
$ws.Range("R6").Value = 16.378560778492
$ws.Range("S6").Value = 0.0009717098825491687
$ws.Range("T6").Value = 0.0009717098825491687
$ws.Range("I7").Value = 0.396815038797359
$ws.Range("J7").Value = 0.396815038797359
$ws.Range("M7").Value = 1.646588666666666
$ws.Range("N7").Value = 4.939766
$ws.Range("O7").Value = 0.039310317935267
$ws.Range("P7").Value = 0.039310317935267
$ws.Range("Q7").Value = 29.21401762367111
$ws.Range("R7").Value = 262.9261586130399
$ws.Range("S7").Value = 0.01559892533661949
$ws.Range("T7").Value = 0.01559892533661949
$ws.Range("I8").Value = 0.396815038797359
$ws.Range("J8").Value = 0.396815038797359
$ws.Range("O8").Value = 0.278787195370394
$ws.Range("P8").Value = 0.278787195370394
$ws.Range("S8").Value = 0.1106269517471098
$ws.Range("T8").Value = 0.1106269517471098
$ws.Range("I9").Value = 0.396815038797359
$ws.Range("J9").Value = 0.396815038797359
$ws.Range("M9").Value = 0.7553226666666667
$ws.Range("N9").Value = 2.265968
$ws.Range("O9").Value = 0.01803241742850595
$ws.Range("P9").Value = 0.01803241742850595
$ws.Range("Q9").Value = 13.40104553265778
$ws.Range("R9").Value = 120.60940979392
$ws.Range("S9").Value = 0.00715553442150276
$ws.Range("T9").Value = 0.00715553442150276
$ws.Range("I10").Value = 0.396815038797359
$ws.Range("J10").Value = 0.396815038797359
$ws.Range("M10").Value = 27.21325766666666
$ws.Range("N10").Value = 81.63977299999999
$ws.Range("O10").Value = 0.6496836961088899
$ws.Range("P10").Value = 0.6496836961088899
$ws.Range("Q10").Value = 482.8216087997911
$ws.Range("R10").Value = 4345.394479198119
$ws.Range("S10").Value = 0.2578042610774607
$ws.Range("T10").Value = 0.2578042610774607
$ws.Range("I11").Value = 0.396815038797359
$ws.Range("J11").Value = 0.396815038797359
$ws.Range("M11").Value = 0.5942236666666667
$ws.Range("N11").Value = 1.782671
$ws.Range("O11").Value = 0.01418637315694314
$ws.Range("P11").Value = 0.01418637315694314
$ws.Range("Q11").Value = 10.54280344680445
$ws.Range("R11").Value = 94.88523102124
$ws.Range("S11").Value = 0.005629366214666203
$ws.Range("T11").Value = 0.005629366214666203
$ws.Range("G12").Value = 13.27534766666667
$ws.Range("H12").Value = 39.826043
$ws.Range("I12").Value = 0.2969120759943797
$ws.Range("J12").Value = 0.2969120759943796
$ws.Range("M12").Value = 1.646588666666666
$ws.Range("N12").Value = 4.939766
$ws.Range("O12").Value = 0.039310317935267
$ws.Range("P12").Value = 0.039310317935267
$ws.Range("Q12").Value = 21.85903701399311
$ws.Range("R12").Value = 196.731333125938
$ws.Range("S12").Value = 0.01167170810615922
$ws.Range("T12").Value = 0.01167170810615922
$ws.Range("G13").Value = 13.27534766666667
$ws.Range("H13").Value = 39.826043
$ws.Range("I13").Value = 0.2969120759943797
$ws.Range("J13").Value = 0.2969120759943796
$ws.Range("O13").Value = 0.278787195370394
$ws.Range("P13").Value = 0.278787195370394
$ws.Range("Q13").Value = 155.0234122416385
$ws.Range("R13").Value = 1395.210710174746
$ws.Range("S13").Value = 0.0827752849380744
$ws.Range("T13").Value = 0.08277528493807437
$ws.Range("G14").Value = 13.27534766666667
$ws.Range("H14").Value = 39.826043
$ws.Range("I14").Value = 0.2969120759943797
$ws.Range("J14").Value = 0.2969120759943796
$ws.Range("M14").Value = 0.7553226666666667
$ws.Range("N14").Value = 2.265968
$ws.Range("O14").Value = 0.01803241742850595
$ws.Range("P14").Value = 0.01803241742850595
$ws.Range("Q14").Value = 10.02717100051378
$ws.Range("R14").Value = 90.244539004624
$ws.Range("S14").Value = 0.005354042493894935
$ws.Range("T14").Value = 0.005354042493894933
$ws.Range("G15").Value = 13.27534766666667
$ws.Range("H15").Value = 39.826043
$ws.Range("I15").Value = 0.2969120759943797
$ws.Range("J15").Value = 0.2969120759943796
$ws.Range("M15").Value = 27.21325766666666
$ws.Range("N15").Value = 81.63977299999999
$ws.Range("O15").Value = 0.6496836961088899
$ws.Range("P15").Value = 0.6496836961088899
$ws.Range("Q15").Value = 361.2654566675821
$ws.Range("R15").Value = 3251.389110008239
$ws.Range("S15").Value = 0.1928989349513922
$ws.Range("T15").Value = 0.1928989349513922
$ws.Range("G16").Value = 13.27534766666667
$ws.Range("H16").Value = 39.826043
$ws.Range("I16").Value = 0.2969120759943797
$ws.Range("J16").Value = 0.2969120759943796
$ws.Range("M16").Value = 0.5942236666666667
$ws.Range("N16").Value = 1.782671
$ws.Range("O16").Value = 0.01418637315694314
$ws.Range("P16").Value = 0.01418637315694314
$ws.Range("Q16").Value = 7.888525766761445
$ws.Range("R16").Value = 70.996731900853
$ws.Range("S16").Value = 0.004212105504858929
$ws.Range("T16").Value = 0.004212105504858928
$ws.Range("G17").Value = 3.455866
$ws.Range("H17").Value = 10.367598
$ws.Range("I17").Value = 0.07729276657626213
$ws.Range("J17").Value = 0.07729276657626213
$ws.Range("M17").Value = 1.646588666666666
$ws.Range("N17").Value = 4.939766
$ws.Range("O17").Value = 0.039310317935267
$ws.Range("P17").Value = 0.039310317935267
$ws.Range("Q17").Value = 5.690389789118666
$ws.Range("R17").Value = 51.21350810206799
$ws.Range("S17").Value = 0.003038403228209243
$ws.Range("T17").Value = 0.003038403228209243
$ws.Range("G18").Value = 3.455866
$ws.Range("H18").Value = 10.367598
$ws.Range("I18").Value = 0.07729276657626213
$ws.Range("J18").Value = 0.07729276657626213
$ws.Range("O18").Value = 0.278787195370394
$ws.Range("P18").Value = 0.278787195370394
$ws.Range("Q18").Value = 40.35601575355066
$ws.Range("R18").Value = 363.204141781956
$ws.Range("S18").Value = 0.02154823361621465
$ws.Range("T18").Value = 0.02154823361621465
$ws.Range("G19").Value = 3.455866
$ws.Range("H19").Value = 10.367598
$ws.Range("I19").Value = 0.07729276657626213
$ws.Range("J19").Value = 0.07729276657626213
$ws.Range("M19").Value = 0.7553226666666667
$ws.Range("N19").Value = 2.265968
$ws.Range("O19").Value = 0.01803241742850595
$ws.Range("P19").Value = 0.01803241742850595
$ws.Range("Q19").Value = 2.610293922762667
$ws.Range("R19").Value = 23.492645304864
$ws.Range("S19").Value = 0.001393775431107231
$ws.Range("T19").Value = 0.001393775431107231
$ws.Range("G20").Value = 3.455866
$ws.Range("H20").Value = 10.367598
$ws.Range("I20").Value = 0.07729276657626213
$ws.Range("J20").Value = 0.07729276657626213
$ws.Range("M20").Value = 27.21325766666666
$ws.Range("N20").Value = 81.63977299999999
$ws.Range("O20").Value = 0.6496836961088899
$ws.Range("P20").Value = 0.6496836961088899
$ws.Range("Q20").Value = 94.04537191947266
$ws.Range("R20").Value = 846.4083472752538
$ws.Range("S20").Value = 0.05021585027174765
$ws.Range("T20").Value = 0.05021585027174765
$ws.Range("G21").Value = 3.455866
$ws.Range("H21").Value = 10.367598
$ws.Range("I21").Value = 0.07729276657626213
$ws.Range("J21").Value = 0.07729276657626213
$ws.Range("M21").Value = 0.5942236666666667
$ws.Range("N21").Value = 1.782671
$ws.Range("O21").Value = 0.01418637315694314
$ws.Range("P21").Value = 0.01418637315694314
$ws.Range("Q21").Value = 2.053557366028667
$ws.Range("R21").Value = 18.482016294258
$ws.Range("S21").Value = 0.001096504028983357
$ws.Range("T21").Value = 0.001096504028983357
$ws.Range("G22").Value = 7.175465666666668
$ws.Range("H22").Value = 21.526397
$ws.Range("I22").Value = 0.1604841139238761
$ws.Range("J22").Value = 0.1604841139238761
$ws.Range("M22").Value = 1.646588666666666
$ws.Range("N22").Value = 4.939766
$ws.Range("O22").Value = 0.039310317935267
$ws.Range("P22").Value = 0.039310317935267
$ws.Range("Q22").Value = 11.81504044478911
$ws.Range("R22").Value = 106.335364003102
$ws.Range("S22").Value = 0.006308681541907178
$ws.Range("T22").Value = 0.006308681541907178
$ws.Range("G23").Value = 7.175465666666668
$ws.Range("H23").Value = 21.526397
$ws.Range("I23").Value = 0.1604841139238761
$ws.Range("J23").Value = 0.1604841139238761
$ws.Range("O23").Value = 0.278787195370394
$ws.Range("P23").Value = 0.278787195370394
$ws.Range("Q23").Value = 83.79179212477045
$ws.Range("R23").Value = 754.1261291229342
$ws.Range("S23").Value = 0.0447409160223402
$ws.Range("T23").Value = 0.0447409160223402
$ws.Range("G24").Value = 7.175465666666668
$ws.Range("H24").Value = 21.526397
$ws.Range("I24").Value = 0.1604841139238761
$ws.Range("J24").Value = 0.1604841139238761
$ws.Range("M24").Value = 0.7553226666666667
$ws.Range("N24").Value = 2.265968
$ws.Range("O24").Value = 0.01803241742850595
$ws.Range("P24").Value = 0.01803241742850595
$ws.Range("Q24").Value = 5.419791861921778
$ws.Range("R24").Value = 48.77812675729601
$ws.Range("S24").Value = 0.002893916532919237
$ws.Range("T24").Value = 0.002893916532919237
$ws.Range("G25").Value = 7.175465666666668
$ws.Range("H25").Value = 21.526397
$ws.Range("I25").Value = 0.1604841139238761
$ws.Range("J25").Value = 0.1604841139238761
$ws.Range("M25").Value = 27.21325766666666
$ws.Range("N25").Value = 81.63977299999999
$ws.Range("O25").Value = 0.6496836961088899
$ws.Range("P25").Value = 0.6496836961088899
$ws.Range("Q25").Value = 195.2677960653201
$ws.Range("R25").Value = 1757.410164587881
$ws.Range("S25").Value = 0.104263912300824
$ws.Range("T25").Value = 0.104263912300824
$ws.Range("G26").Value = 7.175465666666668
$ws.Range("H26").Value = 21.526397
$ws.Range("I26").Value = 0.1604841139238761
$ws.Range("J26").Value = 0.1604841139238761
$ws.Range("M26").Value = 0.5942236666666667
$ws.Range("N26").Value = 1.782671
$ws.Range("O26").Value = 0.01418637315694314
$ws.Range("P26").Value = 0.01418637315694314
$ws.Range("Q26").Value = 4.263831518487446
$ws.Range("R26").Value = 38.37448366638701
$ws.Range("S26").Value = 0.00227668752588548
$ws.Range("T26").Value = 0.00227668752588548
